# correção nos dados e inicio da analise PNAD 2009
#
# The sheet had two "section header" rows (row 5: "situação do domicílio",
# row 8: "grandes regiões e unidades da federação") that carried a label in
# column A but no data in B:G. This edit removes those two empty header
# rows (their data rows shift up to take their place), and fixes the
# column header in row 2 so that "total" (not the stray
# "unnamed: 1_level_1" label) is used for column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "situação do domicílio" section-header row (row 5) - it has
# a label but no data; everything below shifts up by one row.
$ws.Rows(5).Delete()

# After the first deletion, the "grandes regiões e unidades da federação"
# section-header row (originally row 8) is now at row 7. Remove it too.
$ws.Rows(7).Delete()

# Fix row 2's first data-column header: it should read "total" to match
# column B (previously it referenced the stray "unnamed: 1_level_1" text).
$ws.Range("B2").Value = "total"
